# Auto-generated script applying numeric corrections to the Cactuar_Profits
# leve-profit tracking workbook (commit: "chore: update Sheets via scheduled runner").
# Each worksheet corresponds to a crafting/gathering job; columns H-N hold
# market-price/profit figures that the scheduled runner refreshed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 159165.22
$ws.Range("I132").Value = 349171.28
$ws.Range("K132").Value = 1047513.84
$ws.Range("M132").Value = -1044983.84
$ws.Range("H138").Value = 6870.027
$ws.Range("J138").Value = 8679.440000000001
$ws.Range("L138").Value = 26038.32
$ws.Range("N138").Value = -36318.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1994.6
$ws.Range("I45").Value = 1944.5
$ws.Range("J45").Value = 2195
$ws.Range("K45").Value = 1944.5
$ws.Range("L45").Value = 2195
$ws.Range("M45").Value = -1567.5
$ws.Range("N45").Value = -2949
$ws.Range("H132").Value = 4870
$ws.Range("I132").Value = 2174.913
$ws.Range("K132").Value = 6524.739
$ws.Range("M132").Value = -3994.739

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 6946443.5
$ws.Range("I99").Value = 10418166
$ws.Range("K99").Value = 10418166
$ws.Range("M99").Value = -10416668
$ws.Range("H107").Value = 1615.5652
$ws.Range("I107").Value = 2278.5
$ws.Range("J107").Value = 1105.6154
$ws.Range("K107").Value = 2278.5
$ws.Range("L107").Value = 1105.6154
$ws.Range("M107").Value = -358.5
$ws.Range("N107").Value = -4945.6154
$ws.Range("H134").Value = 3573.3809
$ws.Range("I134").Value = 3088.0667
$ws.Range("K134").Value = 9264.2001
$ws.Range("M134").Value = -6729.2001
$ws.Range("H140").Value = 151803
$ws.Range("J140").Value = 151803
$ws.Range("L140").Value = 151803
$ws.Range("N140").Value = -162163

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3755.2646
$ws.Range("I31").Value = 3199.9607
$ws.Range("J31").Value = 5421.1763
$ws.Range("K31").Value = 3199.9607
$ws.Range("L31").Value = 5421.1763
$ws.Range("M31").Value = -2904.9607
$ws.Range("N31").Value = -6011.1763
$ws.Range("H34").Value = 3755.2646
$ws.Range("I34").Value = 3199.9607
$ws.Range("J34").Value = 5421.1763
$ws.Range("K34").Value = 3199.9607
$ws.Range("L34").Value = 5421.1763
$ws.Range("M34").Value = -2997.9607
$ws.Range("N34").Value = -5825.1763
$ws.Range("H58").Value = 2186.2222
$ws.Range("I58").Value = 2216.3333
$ws.Range("J58").Value = 2126
$ws.Range("K58").Value = 2216.3333
$ws.Range("L58").Value = 2126
$ws.Range("M58").Value = -2013.3333
$ws.Range("N58").Value = -2532
$ws.Range("H62").Value = 101250
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 133333.33
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 133333.33
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -134581.33
$ws.Range("H65").Value = 101250
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 133333.33
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 666666.6499999999
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -672906.6499999999
$ws.Range("H92").Value = 75000
$ws.Range("J92").Value = 75000
$ws.Range("L92").Value = 75000
$ws.Range("N92").Value = -79992
$ws.Range("H136").Value = 2186.2222
$ws.Range("I136").Value = 2216.3333
$ws.Range("J136").Value = 2126
$ws.Range("K136").Value = 6648.999899999999
$ws.Range("L136").Value = 6378
$ws.Range("M136").Value = -4098.999899999999
$ws.Range("N136").Value = -11478

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1014.1429
$ws.Range("I92").Value = 883.1667
$ws.Range("J92").Value = 1800
$ws.Range("K92").Value = 2649.5001
$ws.Range("L92").Value = 5400
$ws.Range("M92").Value = -1401.5001
$ws.Range("N92").Value = -7896
$ws.Range("H96").Value = 12250
$ws.Range("I96").Value = 10000
$ws.Range("J96").Value = 13000
$ws.Range("K96").Value = 30000
$ws.Range("L96").Value = 39000
$ws.Range("M96").Value = -27941
$ws.Range("N96").Value = -43118
$ws.Range("H98").Value = 1899.8572
$ws.Range("J98").Value = 1999.8334
$ws.Range("L98").Value = 5999.5002
$ws.Range("N98").Value = -8995.5002
$ws.Range("H114").Value = 1952.7778
$ws.Range("J114").Value = 1952.7778
$ws.Range("L114").Value = 5858.3334
$ws.Range("N114").Value = -12366.3334
$ws.Range("H122").Value = 468.8125
$ws.Range("I122").Value = 423.5
$ws.Range("K122").Value = 3811.5
$ws.Range("M122").Value = -1361.5
$ws.Range("H128").Value = 349989
$ws.Range("I128").Value = 349989
$ws.Range("K128").Value = 1049967
$ws.Range("M128").Value = -1044987
$ws.Range("H133").Value = 16897.889
$ws.Range("J133").Value = 15709.322
$ws.Range("L133").Value = 47127.966
$ws.Range("N133").Value = -57247.966
$ws.Range("H137").Value = 6122167.5
$ws.Range("I137").Value = 3584.6
$ws.Range("J137").Value = 8475469
$ws.Range("K137").Value = 10753.8
$ws.Range("L137").Value = 25426407
$ws.Range("M137").Value = -5653.799999999999
$ws.Range("N137").Value = -25436607

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 22498.25
$ws.Range("I58").Value = 19997.666
$ws.Range("K58").Value = 19997.666
$ws.Range("M58").Value = -19720.666
$ws.Range("H80").Value = 1399775.4
$ws.Range("I80").Value = 2380629.5
$ws.Range("J80").Value = 26579.4
$ws.Range("K80").Value = 2380629.5
$ws.Range("L80").Value = 26579.4
$ws.Range("M80").Value = -2379631.5
$ws.Range("N80").Value = -28575.4
$ws.Range("H83").Value = 1399775.4
$ws.Range("I83").Value = 2380629.5
$ws.Range("J83").Value = 26579.4
$ws.Range("K83").Value = 11903147.5
$ws.Range("L83").Value = 132897
$ws.Range("M83").Value = -11898155.5
$ws.Range("N83").Value = -142881
$ws.Range("H122").Value = 788401.6
$ws.Range("I122").Value = 919426.9399999999
$ws.Range("K122").Value = 2758280.82
$ws.Range("M122").Value = -2755830.82
$ws.Range("H132").Value = 9997.75
$ws.Range("I132").Value = 9997
$ws.Range("K132").Value = 29991
$ws.Range("M132").Value = -27461

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1156.6666
$ws.Range("I22").Value = 810.4
$ws.Range("J22").Value = 1329.8
$ws.Range("K22").Value = 810.4
$ws.Range("L22").Value = 1329.8
$ws.Range("M22").Value = -515.4
$ws.Range("N22").Value = -1919.8
$ws.Range("H27").Value = 1156.6666
$ws.Range("I27").Value = 810.4
$ws.Range("J27").Value = 1329.8
$ws.Range("K27").Value = 810.4
$ws.Range("L27").Value = 1329.8
$ws.Range("M27").Value = -703.4
$ws.Range("N27").Value = -1543.8
$ws.Range("H40").Value = 3414.524
$ws.Range("I40").Value = 2186.1538
$ws.Range("K40").Value = 2186.1538
$ws.Range("M40").Value = -2050.1538
$ws.Range("H61").Value = 3186.8125
$ws.Range("I61").Value = 3531.6365
$ws.Range("K61").Value = 3531.6365
$ws.Range("M61").Value = -3329.6365
$ws.Range("H113").Value = 3186.8125
$ws.Range("I113").Value = 3531.6365
$ws.Range("K113").Value = 3531.6365
$ws.Range("M113").Value = -1361.6365
$ws.Range("H122").Value = 7481.227
$ws.Range("I122").Value = 3715.5
$ws.Range("K122").Value = 11146.5
$ws.Range("M122").Value = -8696.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 27000
$ws.Range("J59").Value = 27000
$ws.Range("L59").Value = 27000
$ws.Range("N59").Value = -28476
$ws.Range("H62").Value = 4939.8667
$ws.Range("I62").Value = 4058.4167
$ws.Range("K62").Value = 4058.4167
$ws.Range("M62").Value = -3434.4167
$ws.Range("H65").Value = 4939.8667
$ws.Range("I65").Value = 4058.4167
$ws.Range("K65").Value = 20292.0835
$ws.Range("M65").Value = -17172.0835
$ws.Range("H81").Value = 5214303.5
$ws.Range("I81").Value = 4172085.8
$ws.Range("J81").Value = 6951332.5
$ws.Range("K81").Value = 8344171.6
$ws.Range("L81").Value = 13902665
$ws.Range("M81").Value = -8343110.6
$ws.Range("N81").Value = -13904787
$ws.Range("H84").Value = 5214303.5
$ws.Range("I84").Value = 4172085.8
$ws.Range("J84").Value = 6951332.5
$ws.Range("K84").Value = 41720858
$ws.Range("L84").Value = 69513325
$ws.Range("M84").Value = -41715554
$ws.Range("N84").Value = -69523933
$ws.Range("H126").Value = 1430.421
$ws.Range("I126").Value = 1221.6364
$ws.Range("K126").Value = 3664.9092
$ws.Range("M126").Value = -1194.9092
